# "Generate Report for Handoff" — file b.md has been handed off again:
# status flips from "Handed back: in sync with en-US" to "Ready for handoff",
# a new handoff xlf + timestamp is recorded, and an error detail explaining
# the stale handback version is attached, for both the zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

$status      = "Ready for handoff"
$hoDateOverview = "2016-08-25 14:37:51"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fc26fadb063d9f052b50a22571eba7e399e73cdf/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/df7cf2e768d59e80b65915cbf88b879bbdbafeb2/e2e/b.md."

# ---------------------------------------------------------------------------
# Overview sheet: row 3 is b.md
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $status
$overview.Range("F3").Value = $status
$overview.Range("G3").Value = $hoDateOverview

# ---------------------------------------------------------------------------
# zh-cn sheet: row 3 is b.md
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $status
# O3 already holds the text "False" (not a real boolean); copy it so F3 stays
# a text cell too instead of being auto-typed as a COM boolean.
$zhcn.Range("O3").Copy($zhcn.Range("F3"))
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-25 14:37:46"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.16

# ---------------------------------------------------------------------------
# de-de sheet: row 3 is b.md
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $status
# O3 already holds the text "False" (not a real boolean); copy it so F3 stays
# a text cell too instead of being auto-typed as a COM boolean.
$dede.Range("O3").Copy($dede.Range("F3"))
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = $hoDateOverview
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.16
